# Timesheet: add two new clock-in/out entries (rows 20-21) for 2026-02-06
# and move the "Total Duration:" summary row down to row 22 with the
# updated total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage literal text so that date/time-looking
# strings ("2026-02-06", "16:57:32", ...) are not auto-converted into
# date/time serial numbers when written into the sheet.
$scratch = $ws.Cells.Item(100, 26)

function Set-TextCell($targetAddr, $styleSourceAddr, $text) {
    # Stage the literal text as a formula string result (avoids the
    # engine's "looks like a date" auto-detection entirely).
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'

    # Bring over the formatting (style) from an existing, already
    # correctly-styled cell.
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4122)

    # Bring over only the literal value (no formula, no format change).
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
}

# --- Row 20: new entry ---
Set-TextCell "A20" "A19" "2026-02-06"
Set-TextCell "B20" "B19" "16:57:32"
Set-TextCell "C20" "C19" "18:11:58"
Set-TextCell "D20" "D19" "1.24 Hours"

# --- Row 21: new entry ---
Set-TextCell "A21" "A19" "2026-02-06"
Set-TextCell "B21" "B19" "19:15:30"
Set-TextCell "C21" "C19" "19:43:57"
Set-TextCell "D21" "D19" "0.47 Hours"

# --- Row 22: Total Duration summary (moved down from row 20) ---
Set-TextCell "C22" "C20" "Total Duration:"
Set-TextCell "D22" "D20" "28 Hours"

# Clean up the scratch cell completely (value + formatting) so it
# leaves no trace (and doesn't register a new, unused style).
$scratch.Clear()
